# Updates cryptos list figures (prices / 1h volume %) per the Wed Feb 21 2024
# GitHub Actions data refresh, including the Filecoin/Hedera swap at rows 29-30,
# the RocketPoolETH -> TheGraph -> BEAM -> SEI shift at rows 49-51, and all
# updated Price / Volume(1h) cell text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = '51.805.39'
$ws.Cells.Item(2, 5).Value = '  -0.38%  '
$ws.Cells.Item(3, 4).Value = '2.962.66'
$ws.Cells.Item(3, 5).Value = '  +0.95%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '353.58'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.20%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '106.05'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -4.20%  '
$ws.Cells.Item(7, 5).Value = '  -3.07%  '
$ws.Cells.Item(8, 5).Value = '  +0.03%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.600'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -5.50%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '37.58'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -4.70%  '
$ws.Cells.Item(11, 5).Value = '  +2.33%  '
$ws.Cells.Item(12, 5).Value = '  -3.77%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.88'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -4.47%  '
$ws.Cells.Item(14, 4).Value = '3.435.29'
$ws.Cells.Item(14, 5).Value = '  +1.21%  '
$ws.Cells.Item(15, 5).Value = '  -6.31%  '
$ws.Cells.Item(16, 4).Value = '2.946.91'
$ws.Cells.Item(16, 5).Value = '  +0.58%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.984'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.42%  '
$ws.Cells.Item(18, 4).Value = '51.749.27'
$ws.Cells.Item(18, 5).Value = '  -0.55%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.32'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.75%  '
$ws.Cells.Item(20, 5).Value = '  -3.43%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.32'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -5.18%  '
$ws.Cells.Item(22, 4).Value = '0.0₃0960'
$ws.Cells.Item(22, 5).Value = '  -2.50%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '68.92'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -3.05%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '264.38'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -2.36%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.68'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -5.33%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.175'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -6.65%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '26.62'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -1.99%  '
$ws.Cells.Item(28, 5).Value = '  +0.00%  '
$ws.Cells.Item(29, 2).Value = 'Hedera'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.110'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +3.71%  '
$ws.Cells.Item(30, 2).Value = 'Filecoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.26'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -3.78%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.26'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +2.80%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.05'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -5.81%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.17'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -5.25%  '
$ws.Cells.Item(34, 5).Value = '  -7.15%  '
$ws.Cells.Item(35, 5).Value = '  -2.95%  '
$ws.Cells.Item(36, 5).Value = '  -4.41%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -0.10%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.23'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -0.92%  '
$ws.Cells.Item(39, 5).Value = '  +1.41%  '
$ws.Cells.Item(40, 5).Value = '  -5.03%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '17.30'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -6.41%  '
$ws.Cells.Item(42, 5).Value = '  -3.68%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '22.81'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -1.18%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '123.49'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +3.64%  '
$ws.Cells.Item(45, 5).Value = '  -0.30%  '
$ws.Cells.Item(46, 4).Value = '2.105.85'
$ws.Cells.Item(46, 5).Value = '  -1.65%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.26'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -6.46%  '
$ws.Cells.Item(48, 5).Value = '  -7.85%  '
$ws.Cells.Item(49, 2).Value = 'TheGraph'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.237'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -3.94%  '
$ws.Cells.Item(50, 2).Value = 'BEAM'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0325'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -2.56%  '
$ws.Cells.Item(51, 2).Value = 'SEI'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.889'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -2.42%  '
